# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 527 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45175
